$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'60.982.21"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = "'  -0.07%  "
$ws.Range('E2').ClearFormats()
$ws.Range('D3').Value = "'2.920.20"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = "'  -0.13%  "
$ws.Range('E3').ClearFormats()
$ws.Range('E4').Value = "'  -0.01%  "
$ws.Range('E4').ClearFormats()
$ws.Range('D5').Value = "'590.44"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = "'  +0.56%  "
$ws.Range('E5').ClearFormats()
$ws.Range('D6').Value = "'146.62"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = "'  +0.30%  "
$ws.Range('E6').ClearFormats()
$ws.Range('E7').Value = "'  +0.05%  "
$ws.Range('E7').ClearFormats()
$ws.Range('D8').Value = "'0.506"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = "'  +0.14%  "
$ws.Range('E8').ClearFormats()
$ws.Range('D9').Value = "'6.93"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = "'  +0.85%  "
$ws.Range('E9').ClearFormats()
$ws.Range('E10').Value = "'  -1.09%  "
$ws.Range('E10').ClearFormats()
$ws.Range('E11').Value = "'  -1.75%  "
$ws.Range('E11').ClearFormats()
$ws.Range('E12').Value = "'  -0.16%  "
$ws.Range('E12').ClearFormats()
$ws.Range('D13').Value = "'33.62"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = "'  -0.19%  "
$ws.Range('E13').ClearFormats()
$ws.Range('E14').Value = "'  -0.22%  "
$ws.Range('E14').ClearFormats()
$ws.Range('D15').Value = "'3.405.75"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = "'  -0.04%  "
$ws.Range('E15').ClearFormats()
$ws.Range('D16').Value = "'60.966.28"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = "'  -0.01%  "
$ws.Range('E16').ClearFormats()
$ws.Range('E17').Value = "'  -1.35%  "
$ws.Range('E17').ClearFormats()
$ws.Range('D18').Value = "'2.923.75"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = "'  +0.01%  "
$ws.Range('E18').ClearFormats()
$ws.Range('D19').Value = "'432.54"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = "'  +0.39%  "
$ws.Range('E19').ClearFormats()
$ws.Range('D20').Value = "'13.41"
$ws.Range('D20').ClearFormats()
$ws.Range('D21').Value = "'0.677"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = "'  -0.76%  "
$ws.Range('E21').ClearFormats()
$ws.Range('D22').Value = "'7.11"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = "'  -0.43%  "
$ws.Range('E22').ClearFormats()
$ws.Range('E23').Value = "'  +0.92%  "
$ws.Range('E23').ClearFormats()
$ws.Range('D24').Value = "'10.93"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = "'  +0.75%  "
$ws.Range('E24').ClearFormats()
$ws.Range('D25').Value = "'2.21"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = "'  -0.83%  "
$ws.Range('E25').ClearFormats()
$ws.Range('D26').Value = "'11.86"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = "'  -0.83%  "
$ws.Range('E26').ClearFormats()
$ws.Range('E27').Value = "'  -0.06%  "
$ws.Range('E27').ClearFormats()
$ws.Range('D28').Value = "'2.27"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = "'  +4.55%  "
$ws.Range('E28').ClearFormats()
$ws.Range('D29').Value = "'2.60"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = "'  -0.59%  "
$ws.Range('E29').ClearFormats()
$ws.Range('D30').Value = "'6.99"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = "'  -3.19%  "
$ws.Range('E30').ClearFormats()
$ws.Range('B31').Value = "'Hedera"
$ws.Range('B31').ClearFormats()
$ws.Range('C31').Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range('C31').ClearFormats()
$ws.Range('D31').Value = "'0.110"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = "'  +2.80%  "
$ws.Range('E31').ClearFormats()
$ws.Range('B32').Value = "'EthereumClassic"
$ws.Range('B32').ClearFormats()
$ws.Range('C32').Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range('C32').ClearFormats()
$ws.Range('D32').Value = "'26.67"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = "'  +0.25%  "
$ws.Range('E32').ClearFormats()
$ws.Range('D34').Value = "'0.0₃0869"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = "'  -0.63%  "
$ws.Range('E34').ClearFormats()
$ws.Range('E35').Value = "'  -0.32%  "
$ws.Range('E35').ClearFormats()
$ws.Range('E36').Value = "'  -0.43%  "
$ws.Range('E36').ClearFormats()
$ws.Range('D37').Value = "'3.02"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = "'  -0.72%  "
$ws.Range('E37').ClearFormats()
$ws.Range('E38').Value = "'  -1.39%  "
$ws.Range('E38').ClearFormats()
$ws.Range('D39').Value = "'0.121"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = "'  -5.31%  "
$ws.Range('E39').ClearFormats()
$ws.Range('D40').Value = "'8.54"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = "'  -1.67%  "
$ws.Range('E40').ClearFormats()
$ws.Range('D41').Value = "'41.70"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = "'  +0.54%  "
$ws.Range('E41').ClearFormats()
$ws.Range('E42').Value = "'  -5.11%  "
$ws.Range('E42').ClearFormats()
$ws.Range('D43').Value = "'377.26"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = "'  -0.45%  "
$ws.Range('E43').ClearFormats()
$ws.Range('B44').Value = "'Maker"
$ws.Range('B44').ClearFormats()
$ws.Range('C44').Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range('C44').ClearFormats()
$ws.Range('D44').Value = "'2.707.07"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = "'  +0.15%  "
$ws.Range('E44').ClearFormats()
$ws.Range('B45').Value = "'VeChain"
$ws.Range('B45').ClearFormats()
$ws.Range('C45').Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range('C45').ClearFormats()
$ws.Range('D45').Value = "'0.0344"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = "'  -1.71%  "
$ws.Range('E45').ClearFormats()
$ws.Range('D46').Value = "'133.85"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = "'  +0.92%  "
$ws.Range('E46').ClearFormats()
$ws.Range('E47').Value = "'  +0.00%  "
$ws.Range('E47').ClearFormats()
$ws.Range('E48').Value = "'  -4.25%  "
$ws.Range('E48').ClearFormats()
$ws.Range('E49').Value = "'  -0.66%  "
$ws.Range('E49').ClearFormats()
$ws.Range('E50').Value = "'  -3.10%  "
$ws.Range('E50').ClearFormats()
$ws.Range('E51').Value = "'  -0.96%  "
$ws.Range('E51').ClearFormats()
